$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2503
$ws.Range("I58").Value = 257.5
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 772.5
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -622.5
$ws.Range("N58").Value = -12300

$ws.Range("H100").Value = 11906189
$ws.Range("I100").Value = 18519532
$ws.Range("J100").Value = 2169.2
$ws.Range("K100").Value = 18519532
$ws.Range("L100").Value = 2169.2
$ws.Range("M100").Value = -18518991
$ws.Range("N100").Value = -3251.2

$ws.Range("H132").Value = 1892.4333
$ws.Range("I132").Value = 1963.3214
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 5889.9642
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -3359.9642
$ws.Range("N132").Value = -7760

$ws.Range("H135").Value = 2118.0378
$ws.Range("I135").Value = 2109.4211
$ws.Range("J135").Value = 2139.8667
$ws.Range("K135").Value = 18984.7899
$ws.Range("L135").Value = 19258.8003
$ws.Range("M135").Value = -16449.7899
$ws.Range("N135").Value = -24328.8003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H45").Value = 17751.166
$ws.Range("I45").Value = 17751.166
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 17751.166
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -17374.166
$ws.Range("N45").ClearContents()

$ws.Range("H74").Value = 1962.3235
$ws.Range("I74").Value = 1789.1333
$ws.Range("J74").Value = 2099.0527
$ws.Range("K74").Value = 1789.1333
$ws.Range("L74").Value = 2099.0527
$ws.Range("M74").Value = -915.1333
$ws.Range("N74").Value = -3847.0527

$ws.Range("H77").Value = 1962.3235
$ws.Range("I77").Value = 1789.1333
$ws.Range("J77").Value = 2099.0527
$ws.Range("K77").Value = 8945.666499999999
$ws.Range("L77").Value = 10495.2635
$ws.Range("M77").Value = -4577.666499999999
$ws.Range("N77").Value = -19231.2635

$ws.Range("H132").Value = 3286.8206
$ws.Range("I132").Value = 1746.0435
$ws.Range("K132").Value = 5238.1305
$ws.Range("M132").Value = -2708.1305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1123.5
$ws.Range("I34").Value = 966.6667
$ws.Range("J34").Value = 1175.7778
$ws.Range("K34").Value = 2900.0001
$ws.Range("L34").Value = 3527.3334
$ws.Range("M34").Value = -2816.0001
$ws.Range("N34").Value = -3695.3334

$ws.Range("H39").Value = 1100
$ws.Range("J39").Value = 1100
$ws.Range("L39").Value = 3300
$ws.Range("N39").Value = -3888

$ws.Range("H68").Value = 2306.5
$ws.Range("I68").Value = 683.3333
$ws.Range("J68").Value = 3280.4
$ws.Range("K68").Value = 2049.9999
$ws.Range("L68").Value = 9841.200000000001
$ws.Range("M68").Value = -1238.9999
$ws.Range("N68").Value = -11463.2

$ws.Range("H71").Value = 2306.5
$ws.Range("I71").Value = 683.3333
$ws.Range("J71").Value = 3280.4
$ws.Range("K71").Value = 6149.9997
$ws.Range("L71").Value = 29523.6
$ws.Range("M71").Value = -2093.9997
$ws.Range("N71").Value = -37635.60000000001

$ws.Range("H131").Value = 1786739.1
$ws.Range("I131").Value = 10000840
$ws.Range("K131").Value = 30002520
$ws.Range("M131").Value = -29997480

$ws.Range("H132").Value = 4500
$ws.Range("J132").Value = 3333.3333
$ws.Range("L132").Value = 29999.9997
$ws.Range("N132").Value = -35059.9997

$ws.Range("H134").Value = 9565.925999999999
$ws.Range("I134").Value = 9162.857
$ws.Range("K134").Value = 27488.571
$ws.Range("M134").Value = -22418.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1659.375
$ws.Range("I43").Value = 682.1429000000001
$ws.Range("K43").Value = 682.1429000000001
$ws.Range("M43").Value = -531.1429000000001

$ws.Range("H46").Value = 12520.5
$ws.Range("I46").Value = 10041
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 10041
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -9885
$ws.Range("N46").Value = -15312

$ws.Range("H57").Value = 8056.75
$ws.Range("J57").Value = 8056.75
$ws.Range("L57").Value = 8056.75
$ws.Range("N57").Value = -9696.75

$ws.Range("H70").Value = 6858.9585
$ws.Range("I70").Value = 7089.05
$ws.Range("J70").Value = 5708.5
$ws.Range("K70").Value = 7089.05
$ws.Range("L70").Value = 5708.5
$ws.Range("M70").Value = -6819.05
$ws.Range("N70").Value = -6248.5

$ws.Range("H73").Value = 6858.9585
$ws.Range("I73").Value = 7089.05
$ws.Range("J73").Value = 5708.5
$ws.Range("K73").Value = 7089.05
$ws.Range("L73").Value = 5708.5
$ws.Range("M73").Value = -6153.05
$ws.Range("N73").Value = -7580.5

$ws.Range("H80").Value = 2480.8696
$ws.Range("I80").Value = 2468.5715
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 2468.5715
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -1470.5715
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 2480.8696
$ws.Range("I83").Value = 2468.5715
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 12342.8575
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -7350.8575
$ws.Range("N83").Value = -22484

$ws.Range("H113").Value = 47620650
$ws.Range("I113").Value = 83334360
$ws.Range("J113").Value = 2364.4443
$ws.Range("K113").Value = 83334360
$ws.Range("L113").Value = 2364.4443
$ws.Range("M113").Value = -83332190
$ws.Range("N113").Value = -6704.4443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 6999.6665
$ws.Range("I18").Value = 6999
$ws.Range("K18").Value = 6999
$ws.Range("M18").Value = -6827

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H40").Value = 52633730
$ws.Range("I40").Value = 76925370
$ws.Range("J40").Value = 1834.1666
$ws.Range("K40").Value = 76925370
$ws.Range("L40").Value = 1834.1666
$ws.Range("M40").Value = -76925234
$ws.Range("N40").Value = -2106.1666

$ws.Range("H46").Value = 1277.6666
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1312.375
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1312.375
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1688.375

$ws.Range("H122").Value = 4794151
$ws.Range("I122").Value = 6497869.5
$ws.Range("J122").Value = 1670666.6
$ws.Range("K122").Value = 19493608.5
$ws.Range("L122").Value = 5011999.800000001
$ws.Range("M122").Value = -19491158.5
$ws.Range("N122").Value = -5016899.800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 35445.8
$ws.Range("J46").Value = 35445.8
$ws.Range("L46").Value = 35445.8
$ws.Range("N46").Value = -35907.8

$ws.Range("H113").Value = 1404.1904
$ws.Range("I113").Value = 1154.0333
$ws.Range("J113").Value = 2029.5834
$ws.Range("K113").Value = 3462.0999
$ws.Range("L113").Value = 6088.7502
$ws.Range("M113").Value = -1292.0999
$ws.Range("N113").Value = -10428.7502

$ws.Range("H134").Value = 35445.8
$ws.Range("J134").Value = 35445.8
$ws.Range("L134").Value = 106337.4
$ws.Range("N134").Value = -111407.4
